$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("wenden", "none", "none"),
    @("spenden", "dog/dog008.jpg", "dog"),
    @("antun", "dog/dog001.jpg", "dog"),
    @("rufen", "none", "none"),
    @("grüßen", "dog/dog015.jpg", "dog"),
    @("tauschen", "house/house027.jpg", "house"),
    @("enden", "none", "none"),
    @("bitten", "dog/dog020.jpg", "dog"),
    @("bieten", "dog/dog014.jpg", "dog"),
    @("opfern", "none", "none"),
    @("altern", "house/house003.jpg", "house"),
    @("fügen", "dog/dog019.jpg", "dog"),
    @("schalten", "none", "none"),
    @("rühren", "house/house009.jpg", "house"),
    @("schneiden", "house/house015.jpg", "house"),
    @("klagen", "none", "none"),
    @("lächeln", "dog/dog002.jpg", "dog"),
    @("lügen", "house/house022.jpg", "house"),
    @("sparen", "none", "none"),
    @("buchen", "house/house028.jpg", "house"),
    @("schützen", "dog/dog026.jpg", "dog"),
    @("hören", "none", "none"),
    @("husten", "house/house023.jpg", "house"),
    @("schütteln", "house/house014.jpg", "house"),
    @("drohen", "none", "none"),
    @("lassen", "house/house018.jpg", "house"),
    @("heben", "dog/dog003.jpg", "dog"),
    @("orten", "none", "none"),
    @("treten", "house/house031.jpg", "house"),
    @("sammeln", "dog/dog012.jpg", "dog"),
    @("weigern", "none", "none"),
    @("brauchen", "house/house012.jpg", "house"),
    @("küssen", "dog/dog010.jpg", "dog"),
    @("ächzen", "none", "none"),
    @("gelten", "house/house030.jpg", "house"),
    @("regnen", "dog/dog031.jpg", "dog"),
    @("kosten", "none", "none"),
    @("zählen", "house/house016.jpg", "house"),
    @("wachsen", "house/house029.jpg", "house"),
    @("dauern", "none", "none"),
    @("münzen", "dog/dog005.jpg", "dog"),
    @("ruhen", "house/house017.jpg", "house"),
    @("stören", "none", "none"),
    @("schicken", "dog/dog017.jpg", "dog"),
    @("kleben", "dog/dog018.jpg", "dog"),
    @("bremsen", "none", "none"),
    @("knien", "house/house004.jpg", "house"),
    @("mauern", "dog/dog000.jpg", "dog")
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}
